$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a "price"-style cell's value while preserving its original
# text (string) storage type. Cells that already contain separators other
# than a single decimal point (e.g. "44.755.40") are left alone, because
# the engine keeps those as text automatically. Cells whose new value looks
# like a plain decimal number would otherwise get silently re-typed as a
# numeric cell, so we force the Text number format on those first (mirrors
# what Excel itself does when you type a number into a Text-formatted cell).
function Set-TextValue {
    param($range, [string]$text)
    $range.NumberFormat = "@"
    $range.Value = $text
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "44.755.40"
$ws.Range("E2").Value = "  +3.90%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.424.13"
$ws.Range("E3").Value = "  +2.34%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.07%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "316.15"

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "101.78"
$ws.Range("E6").Value = "  +6.94%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.515"
$ws.Range("E7").Value = "  +2.51%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.530"
$ws.Range("E9").Value = "  +10.05%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "35.47"
$ws.Range("E10").Value = "  +3.29%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.0802"
$ws.Range("E11").Value = "  +1.85%  "

# Row 12 - Chainlink
Set-TextValue $ws.Range("D12") "18.78"
$ws.Range("E12").Value = "  +1.04%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  -1.17%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "6.94"

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.801.76"
$ws.Range("E15").Value = "  +2.39%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.470.04"
$ws.Range("E16").Value = "  +4.38%  "

# Row 17 - Polygon
Set-TextValue $ws.Range("D17") "0.833"
$ws.Range("E17").Value = "  +4.47%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "44.602.02"
$ws.Range("E18").Value = "  +3.48%  "

# Row 19 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D19") "12.29"
$ws.Range("E19").Value = "  +2.58%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  +1.43%  "

# Row 21 - ShibaInu
$ws.Range("E21").Value = "  +3.55%  "

# Row 22 - Litecoin
Set-TextValue $ws.Range("D22") "68.78"
$ws.Range("E22").Value = "  +0.94%  "

# Row 23 - BitcoinCash
Set-TextValue $ws.Range("D23") "242.15"
$ws.Range("E23").Value = "  +2.83%  "

# Row 24 - ImmutableX
$ws.Range("E24").Value = "  +4.16%  "

# Row 25 - PancakeSwap
Set-TextValue $ws.Range("D25") "2.50"
$ws.Range("E25").Value = "  +2.13%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.14%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "25.23"
$ws.Range("E27").Value = "  +2.85%  "

# Row 28 - Toncoin
Set-TextValue $ws.Range("D28") "2.29"
$ws.Range("E28").Value = "  -3.47%  "

# Row 29 - Cosmos
Set-TextValue $ws.Range("D29") "9.47"
$ws.Range("E29").Value = "  +1.34%  "

# Row 30 - InjectiveProtocol
Set-TextValue $ws.Range("D30") "33.68"
$ws.Range("E30").Value = "  +4.11%  "

# Row 31 - OKB
Set-TextValue $ws.Range("D31") "48.56"
$ws.Range("E31").Value = "  +1.54%  "

# Row 32 - Kaspa
Set-TextValue $ws.Range("D32") "0.127"
$ws.Range("E32").Value = "  +18.31%  "

# Row 33 - Celestia
Set-TextValue $ws.Range("D33") "19.47"
$ws.Range("E33").Value = "  +11.06%  "

# Row 34 - now Filecoin (was Hedera)
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D34") "5.17"
$ws.Range("E34").Value = "  +2.89%  "

# Row 35 - now Hedera (was Filecoin)
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D35") "0.0773"
$ws.Range("E35").Value = "  +6.28%  "

# Row 37 - ARBITRUM
$ws.Range("E37").Value = "  +3.07%  "

# Row 38 - RenderToken
$ws.Range("E38").Value = "  +3.18%  "

# Row 40 - Monero
Set-TextValue $ws.Range("D40") "122.75"
$ws.Range("E40").Value = "  -3.10%  "

# Row 41 - Stellar
$ws.Range("E41").Value = "  +1.68%  "

# Row 42 - WEMIXToken
Set-TextValue $ws.Range("D42") "2.21"
$ws.Range("E42").Value = "  -2.91%  "

# Row 43 - EnergySwap
Set-TextValue $ws.Range("D43") "20.98"
$ws.Range("E43").Value = "  +0.80%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  +4.30%  "

# Row 45 - Maker
$ws.Range("D45").Value = "1.939.81"
$ws.Range("E45").Value = "  +0.24%  "

# Row 46 - ApeXProtocol
$ws.Range("E46").Value = "  -0.24%  "

# Row 47 - NEARProtocol
$ws.Range("E47").Value = "  +8.58%  "

# Row 48 - FraxShare
Set-TextValue $ws.Range("D48") "9.41"
$ws.Range("E48").Value = "  +1.00%  "

# Row 49 - Stacks
Set-TextValue $ws.Range("D49") "1.73"
$ws.Range("E49").Value = "  +14.38%  "

# Row 50 - BitcoinSV
Set-TextValue $ws.Range("D50") "75.35"
$ws.Range("E50").Value = "  +4.93%  "

# Row 51 - MultiversX
Set-TextValue $ws.Range("D51") "54.07"
$ws.Range("E51").Value = "  +5.46%  "
